# Generate Report for Archive
#
# 1. Update the "Ready for handoff" status text (stored as a shared string,
#    so this single change propagates to every cell that references it:
#    Overview!E2:F4, zh-cn!C2:C4, de-de!C2:C4).
# 2. Narrow the "Status" columns (Overview E:F, zh-cn C, de-de C) from
#    17.2159881591797 to 13.4101845877511.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $found = $used.Find("Ready for handoff")
    if ($found -ne $null) {
        $firstAddress = $found.Address()
        do {
            $found.Value = "In Translation"
            $found = $used.FindNext($found)
        } while (($found -ne $null) -and ($found.Address() -ne $firstAddress))
    }
}

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Columns.Item(5).ColumnWidth = 13.4101845877511
$wsOverview.Columns.Item(6).ColumnWidth = 13.4101845877511

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Columns.Item(3).ColumnWidth = 13.4101845877511

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Columns.Item(3).ColumnWidth = 13.4101845877511
